$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the UK flag code from the informal "uk" to the ISO "gb"
$ws.Range("G4").Value = "gb"

# Update percent values (column B) to fix sorting / demo data
$ws.Range("B2").Value = 60
$ws.Range("B3").Value = 30
$ws.Range("B4").Value = 45
$ws.Range("B5").Value = 17
$ws.Range("B6").Value = 24
$ws.Range("B7").Value = 8
$ws.Range("B8").Value = 12
$ws.Range("B9").Value = 13
$ws.Range("B10").Value = 15
$ws.Range("B11").Value = 12
$ws.Range("B12").Value = 5
$ws.Range("B13").Value = 23
$ws.Range("B14").Value = 12
$ws.Range("B15").Value = 20
$ws.Range("B16").Value = 20
$ws.Range("B17").Value = 6
$ws.Range("B18").Value = 5
$ws.Range("B19").Value = 44
$ws.Range("B20").Value = 35
$ws.Range("B21").Value = 14
$ws.Range("B22").Value = 33
$ws.Range("B23").Value = 2
$ws.Range("B24").Value = 12

# Update the active selection to match the author's final cursor position
$ws.Range("B24").Select()
